# Update header labels on the two existing sheets.
$wb = $excel.ActiveWorkbook

$weekly = $wb.Worksheets.Item("Weekly Quantity")
$weekly.Range("B1").Value = "Weekly_PO_Qty"

$monthly = $wb.Worksheets.Item("Monthly Trend")
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# Add a new "PO Forecast" sheet at the end of the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$forecast = $wb.Worksheets.Add($null, $lastSheet)
$forecast.Name = "PO Forecast"

# Reuse the bold/centered/bordered header style and the date-formatted
# column style from the "Weekly Quantity" sheet so no new style entries
# are introduced.
$weekly.Range("A1:B1").Copy()
$forecast.Range("A1:D1").PasteSpecial(-4122)
$weekly.Range("A2").Copy()
$forecast.Range("A2:A41").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row values.
$forecast.Range("A1").Value = "ds"
$forecast.Range("B1").Value = "PO_Forecast"
$forecast.Range("C1").Value = "yhat_lower"
$forecast.Range("D1").Value = "yhat_upper"

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper).
$data = @(
  @(44934.99999999999, 15, -82.55576487472483, 100.4741474169105),
  @(44941.99999999999, 17, -77.87793396801749, 112.6357665856652),
  @(44948.99999999999, 18, -72.98878879295084, 110.8288075110436),
  @(44976.99999999999, 24, -74.54638936851038, 122.2100664609119),
  @(44997.99999999999, 29, -69.52675394294334, 116.691120601032),
  @(45011.99999999999, 32, -61.3351914853712, 120.0594418075132),
  @(45025.99999999999, 35, -53.05822638077785, 128.8677521141265),
  @(45032.99999999999, 36, -55.0998967051517, 129.1090611084082),
  @(45039.99999999999, 38, -53.67185504735634, 122.1224923066083),
  @(45053.99999999999, 41, -57.84711423064734, 132.0748269490673),
  @(45060.99999999999, 42, -54.19740273828238, 131.8770184225355),
  @(45067.99999999999, 44, -47.25033111498943, 134.6377638857275),
  @(45074.99999999999, 45, -45.48494322900258, 132.5019925701082),
  @(45081.99999999999, 46, -48.34275305730831, 137.7249664007612),
  @(45088.99999999999, 48, -44.7972484446492, 142.9146759832799),
  @(45095.99999999999, 49, -40.74077821776905, 142.5580275705223),
  @(45102.99999999999, 51, -37.5580073684613, 146.2085759257067),
  @(45109.99999999999, 52, -38.49335167500242, 148.4231229668963),
  @(45116.99999999999, 54, -42.49973791833027, 146.5705100448327),
  @(45123.99999999999, 55, -37.03998490123382, 146.2277176188048),
  @(45130.99999999999, 57, -37.46907305122018, 153.5360103794071),
  @(45137.99999999999, 58, -34.06492859063184, 153.3068558095313),
  @(45144.99999999999, 60, -28.26227888231784, 155.8024419284108),
  @(45151.99999999999, 61, -32.47931889419844, 147.7891946083631),
  @(45158.99999999999, 63, -28.79510301241102, 156.3952250318258),
  @(45165.99999999999, 64, -27.11805953162811, 156.9733956714904),
  @(45172.99999999999, 66, -28.73357750045034, 163.218302577838),
  @(45179.99999999999, 67, -26.5840973423892, 157.2277861187597),
  @(45186.99999999999, 69, -16.57223547148068, 158.5358652728869),
  @(45193.99999999999, 70, -26.63477937185708, 168.937244390836),
  @(45200.99999999999, 72, -15.74406874780555, 166.4499120325148),
  @(45207.99999999999, 73, -19.45440920858852, 177.217478372237),
  @(45214.99999999999, 75, -22.50818350724959, 172.4730113596667),
  @(45221.99999999999, 76, -15.17975493087823, 161.8645609148895),
  @(45228.99999999999, 78, -15.4096478994843, 171.3124457059841),
  @(45235.99999999999, 79, -10.31503533371121, 171.706670860359),
  @(45242.99999999999, 81, -16.33025366691784, 174.4084756984107),
  @(45249.99999999999, 82, -8.062255811298769, 174.3091075995058),
  @(45256.99999999999, 84, -6.669504974923952, 176.9398882963863),
  @(45263.99999999999, 85, -17.60126562928431, 180.1255934127778)
)

$row = 2
foreach ($entry in $data) {
  $forecast.Range("A$row").Value = $entry[0]
  $forecast.Range("B$row").Value = $entry[1]
  $forecast.Range("C$row").Value = $entry[2]
  $forecast.Range("D$row").Value = $entry[3]
  $row++
}
